$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.309.13"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "2.612.10"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.53"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.23"
$ws.Range("E6").Value = "  -3.64%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.570"
$ws.Range("E8").Value = "  -4.89%  "
$ws.Range("D9").Value = "2.612.39"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.33"
$ws.Range("E10").Value = "  -4.61%  "
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.342"
$ws.Range("E12").Value = "  -1.36%  "
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("D14").Value = "3.062.99"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("D15").Value = "60.301.05"
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.35"
$ws.Range("E16").Value = "  -1.82%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000139"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").Value = "2.614.85"
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.63"
$ws.Range("E19").Value = "  -2.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "343.91"
$ws.Range("E20").Value = "  -4.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.43"
$ws.Range("E21").Value = "  -1.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.10"
$ws.Range("E22").Value = "  -2.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.992"
$ws.Range("E23").Value = "  -0.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.22"
$ws.Range("E24").Value = "  -1.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.419"
$ws.Range("E25").Value = "  -2.53%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.163"
$ws.Range("E26").Value = "  -1.75%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.994"
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("D28").Value = "0.0₃0811"
$ws.Range("E28").Value = "  -3.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.07"
$ws.Range("E29").Value = "  -3.76%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.01"
$ws.Range("E31").Value = "  -2.11%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.59"
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.95"
$ws.Range("E33").Value = "  -2.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.97"
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.97"
$ws.Range("E35").Value = "  -2.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.922"
$ws.Range("E36").Value = "  -2.00%  "
$ws.Range("E37").Value = "  -4.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.867"
$ws.Range("E38").Value = "  +3.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.52"
$ws.Range("E39").Value = "  +0.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.44"
$ws.Range("E40").Value = "  -3.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.63"
$ws.Range("E41").Value = "  -3.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "288.28"
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.624"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.101"
$ws.Range("E44").Value = "  -1.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.997"
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0548"
$ws.Range("E46").Value = "  -1.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.54"
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("E48").Value = "  -1.13%  "
$ws.Range("E49").Value = "  +0.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.70"
$ws.Range("E50").Value = "  -5.29%  "
$ws.Range("D51").Value = "1.955.40"
$ws.Range("E51").Value = "  -1.37%  "
